$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 86
$ws1.Range("F4").Value = 257
$ws1.Range("F5").Value = 150
$ws1.Range("F6").Value = 245
$ws1.Range("F7").Value = 197
$ws1.Range("F8").Value = 1904
$ws1.Range("F9").Value = 342
$ws1.Range("F10").Value = 4444
$ws1.Range("F11").Value = 66
$ws1.Range("F12").Value = 311

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 12

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 86
$ws4.Range("F6").Value = 257
$ws4.Range("F7").Value = 150
$ws4.Range("F8").Value = 245
$ws4.Range("F9").Value = 197
$ws4.Range("F11").Value = 12
$ws4.Range("F12").Value = 1904
$ws4.Range("F13").Value = 342
$ws4.Range("F14").Value = 4444
$ws4.Range("F15").Value = 66
$ws4.Range("F16").Value = 311
